# The pptx table writer now honors explicit/internal column widths instead
# of always evenly dividing the available width across all columns. For the
# table on slide 6 (the "Content Placeholder 5" table with cells "1"/"2"
# located at the left of the slide) this means its two grid columns grow
# from 2501900 EMU (197 pt) to 2514600 EMU (198 pt) each - matching the
# width already used by the sibling table further right on the same slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(6)

for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shp = $s.Shapes.Item($i)
    if ($shp.HasTable) {
        $tbl = $shp.Table
        # Only the first (left-hand) table needs correcting - its columns
        # are narrower (197pt) than the already-correct second table (198pt).
        if ($tbl.Columns.Item(1).Width -eq 197) {
            for ($c = 1; $c -le $tbl.Columns.Count; $c++) {
                $tbl.Columns.Item($c).Width = 198
            }
        }
    }
}
